# Update countries & provincias Spain
# Applies the COVID data refresh: updates the "last updated" timestamp and
# refreshes the per-country stats. Several countries swapped rank (and thus
# row position) in the refreshed data set, so those rows get both a new
# country name and new stats; others just get refreshed numbers in place.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Timestamp header
$ws.Range("A1").Value = "Datos actualizados a 23 de Marzo de 2020 a las 11:46"

# Row data as: Row -> (Country, Casos totales, Nuevos casos, Casos activos, Recuperados, Casos criticos, Muertes hoy, Muertes)
$rows = @(
    @{ Row = 8;   Country = "Alemania";               B = 26198; C = 1325; D = 266; E = 25821; F = 23;  G = 17; H = 111 },
    @{ Row = 9;   Country = "Iran";                    B = 23049; C = 1411; D = 8376; E = 12861; F = 0;  G = 127; H = 1812 },
    @{ Row = 12;  Country = "Suiza";                   B = 7806;  C = 332;  D = 131; E = 7575;  F = 141; G = 2;  H = 100 },
    @{ Row = 15;  Country = "Belgica";                 B = 3743;  C = 342;  D = 350; E = 3305;  F = 322; G = 13; H = 88 },
    @{ Row = 16;  Country = "Austria";                 B = 3679;  C = 97;   D = 9;   E = 3654;  F = 13;  G = 0;  H = 16 },
    @{ Row = 17;  Country = "Noruega";                 B = 2415;  C = 30;   D = 6;   E = 2401;  F = 32;  G = 1;  H = 8 },
    @{ Row = 18;  Country = "Suecia";                  B = 1934;  C = 0;    D = 16;  E = 1893;  F = 80;  G = 4;  H = 25 },
    @{ Row = 81;  Country = "Republica de Macedonia";  B = 115;   C = 0;    D = 1;   E = 112;   F = 1;   G = 1;  H = 2 },
    @{ Row = 111; Country = "Nigeria";                 B = 36;    C = 6;    D = 2;   E = 33;    F = 0;   G = 1;  H = 1 },
    @{ Row = 112; Country = "Cuba";                    B = 35;    C = 0;    D = 0;   E = 34;    F = 0;   G = 0;  H = 1 },
    @{ Row = 121; Country = "Macao";                   B = 24;    C = 2;    D = 10;  E = 14;    F = 0;   G = 0;  H = 0 },
    @{ Row = 122; Country = "Puerto Rico";              B = 23;    C = 0;    D = 0;   E = 22;    F = 0;   G = 0;  H = 1 },
    @{ Row = 123; Country = "Monaco";                  B = 23;    C = 0;    D = 1;   E = 22;    F = 0;   G = 0;  H = 0 },
    @{ Row = 125; Country = "Montenegro";              B = 22;    C = 1;    D = 0;   E = 21;    F = 0;   G = 1;  H = 1 }
)

foreach ($r in $rows) {
    $rowNum = $r.Row
    $ws.Range("A$rowNum").Value = $r.Country
    $ws.Range("B$rowNum").Value = $r.B
    $ws.Range("C$rowNum").Value = $r.C
    $ws.Range("D$rowNum").Value = $r.D
    $ws.Range("E$rowNum").Value = $r.E
    $ws.Range("F$rowNum").Value = $r.F
    $ws.Range("G$rowNum").Value = $r.G
    $ws.Range("H$rowNum").Value = $r.H
}
